# [CHANGED] resources/users - changed password to 12345 hashed
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

$passwordHash = '$2b$10$x0/2y8nbrC55hqKq3jxuBuMt0C1QBXwoUveetxb6U2kuGAxmdWbQ6'

for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 4).Value = $passwordHash
}

$ws.Range("D13").Select()
